$d = $word.ActiveDocument

# The document contains several "<id>...</id>" tags that were split across
# multiple runs (one run per font/style change while the id text was being
# typed/pasted). Re-downloading the tc/tcn/tl content collapses each of
# these back into a single run using the formatting of the opening "<id>"
# run. Re-run the same Find & Replace for each affected id so Word merges
# the runs and restyles the whole tag with the Courier New "id" formatting.

$d.Content.Find.Execute("<id>p062v_1</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p062v_1</id>", 2)

$d.Content.Find.Execute("<id>p062v_2</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p062v_2</id>", 2)

$d.Content.Find.Execute("<id>p062v_3</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p062v_3</id>", 2)
